# Scheduled-runner data refresh: updates currentAveragePrice /
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ and LeveProfitNQ/HQ
# (columns H-N) on each job sheet with freshly pulled market-board
# averages. Generated from the upstream price snapshot diff.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 100
$ws.Range("H100").Value = 4034.8
$ws.Range("I100").Value = 2392.3333
$ws.Range("K100").Value = 2392.3333
$ws.Range("M100").Value = -1851.3333
# Row 132
$ws.Range("H132").Value = 2676.575
$ws.Range("I132").Value = 2407.1353
$ws.Range("K132").Value = 7221.4059
$ws.Range("M132").Value = -4691.4059
# Row 140
$ws.Range("H140").Value = 98228
$ws.Range("J140").Value = 98228
$ws.Range("L140").Value = 98228
$ws.Range("N140").Value = -108588
# Row 141
$ws.Range("H141").Value = 16439.8
$ws.Range("I141").Value = 6600
$ws.Range("K141").Value = 19800
$ws.Range("M141").Value = -14620

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1079.7307
$ws.Range("I32").Value = 831.5714
$ws.Range("K32").Value = 831.5714
$ws.Range("M32").Value = -544.5714
# Row 38
$ws.Range("H38").Value = 20000
$ws.Range("J38").Value = 20000
$ws.Range("L38").Value = 20000
$ws.Range("N38").Value = -20934
# Row 133
$ws.Range("H133").Value = 164052.28
$ws.Range("J133").Value = 183061
$ws.Range("L133").Value = 183061
$ws.Range("N133").Value = -188121
# Row 135
$ws.Range("H135").Value = 96564
$ws.Range("J135").Value = 96564
$ws.Range("L135").Value = 96564
$ws.Range("N135").Value = -106704

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 42560.855
$ws.Range("J82").Value = 89977.336
$ws.Range("L82").Value = 89977.336
$ws.Range("N82").Value = -90743.336
# Row 85
$ws.Range("H85").Value = 42560.855
$ws.Range("J85").Value = 89977.336
$ws.Range("L85").Value = 89977.336
$ws.Range("N85").Value = -92629.336
# Row 105
$ws.Range("H105").Value = 2813.4
$ws.Range("I105").Value = 2813.4
$ws.Range("K105").Value = 2813.4
$ws.Range("M105").Value = -1066.4
# Row 107
$ws.Range("H107").Value = 1497.5518
$ws.Range("I107").Value = 1497.5518
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1497.5518
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 422.4482
$ws.Range("N107").ClearContents()
# Row 132
$ws.Range("H132").Value = 72999.8
$ws.Range("J132").Value = 72999.8
$ws.Range("L132").Value = 72999.8
$ws.Range("N132").Value = -83119.8
# Row 135
$ws.Range("H135").Value = 54798.332
$ws.Range("J135").Value = 54798.332
$ws.Range("L135").Value = 54798.332
$ws.Range("N135").Value = -64938.332
# Row 137
$ws.Range("H137").Value = 106657.78
$ws.Range("J137").Value = 106657.78
$ws.Range("L137").Value = 106657.78
$ws.Range("N137").Value = -116857.78
# Row 138
$ws.Range("H138").Value = 97960
$ws.Range("J138").Value = 97960
$ws.Range("L138").Value = 97960
$ws.Range("N138").Value = -108240

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2868.6428
$ws.Range("I16").Value = 1110
$ws.Range("J16").Value = 4627.2856
$ws.Range("K16").Value = 1110
$ws.Range("L16").Value = 4627.2856
$ws.Range("M16").Value = -823
$ws.Range("N16").Value = -5201.2856
# Row 51
$ws.Range("H51").Value = 27000
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
# Row 60
$ws.Range("H60").Value = 10595
$ws.Range("I60").Value = 13985
$ws.Range("J60").Value = 8900
$ws.Range("K60").Value = 13985
$ws.Range("L60").Value = 8900
$ws.Range("M60").Value = -13474
$ws.Range("N60").Value = -9922
# Row 61
$ws.Range("H61").Value = 27000
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
# Row 105
$ws.Range("H105").Value = 1543.5714
$ws.Range("J105").Value = 2241.8572
$ws.Range("L105").Value = 2241.8572
$ws.Range("N105").Value = -5735.8572
# Row 113
$ws.Range("H113").Value = 2868.6428
$ws.Range("I113").Value = 1110
$ws.Range("J113").Value = 4627.2856
$ws.Range("K113").Value = 1110
$ws.Range("L113").Value = 4627.2856
$ws.Range("M113").Value = 1060
$ws.Range("N113").Value = -8967.285599999999
# Row 134
$ws.Range("H134").Value = 1635.3617
$ws.Range("I134").Value = 1735.4773
$ws.Range("J134").Value = 167
$ws.Range("K134").Value = 5206.4319
$ws.Range("L134").Value = 501
$ws.Range("M134").Value = -2671.4319
$ws.Range("N134").Value = -5571
# Row 138
$ws.Range("H138").Value = 69994.28999999999
$ws.Range("J138").Value = 69994.28999999999
$ws.Range("L138").Value = 69994.28999999999
$ws.Range("N138").Value = -80274.28999999999
# Row 140
$ws.Range("H140").Value = 90260.336
$ws.Range("J140").Value = 90260.336
$ws.Range("L140").Value = 90260.336
$ws.Range("N140").Value = -100620.336

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 46
$ws.Range("H46").Value = 2929.8
$ws.Range("J46").Value = 2929.8
$ws.Range("L46").Value = 8789.400000000001
$ws.Range("N46").Value = -8971.400000000001
# Row 134
$ws.Range("H134").Value = 11978.632
$ws.Range("I134").Value = 2566.5
$ws.Range("K134").Value = 7699.5
$ws.Range("M134").Value = -2629.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 28
$ws.Range("H28").Value = 5900
$ws.Range("J28").Value = 5900
$ws.Range("L28").Value = 5900
$ws.Range("N28").Value = -6284
# Row 47
$ws.Range("H47").Value = 39999
$ws.Range("J47").Value = 39999
$ws.Range("L47").Value = 39999
$ws.Range("N47").Value = -41135
# Row 48
$ws.Range("H48").Value = 39575.668
$ws.Range("J48").Value = 39575.668
$ws.Range("L48").Value = 39575.668
$ws.Range("N48").Value = -40545.668
# Row 93
$ws.Range("H93").Value = 37499
$ws.Range("J93").Value = 37499
$ws.Range("L93").Value = 37499
$ws.Range("N93").Value = -41243
# Row 123
$ws.Range("H123").Value = 49998.25
$ws.Range("J123").Value = 49998.25
$ws.Range("L123").Value = 49998.25
$ws.Range("N123").Value = -54898.25
# Row 128
$ws.Range("H128").Value = 62999.332
$ws.Range("J128").Value = 62999.332
$ws.Range("L128").Value = 62999.332
$ws.Range("N128").Value = -72959.33199999999
# Row 134
$ws.Range("H134").Value = 50829
$ws.Range("J134").Value = 50829
$ws.Range("L134").Value = 152487
$ws.Range("N134").Value = -157557
# Row 135
$ws.Range("H135").Value = 84999.336
$ws.Range("J135").Value = 84999.336
$ws.Range("L135").Value = 84999.336
$ws.Range("N135").Value = -95139.336

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 3074.25
$ws.Range("I61").Value = 3074.25
$ws.Range("K61").Value = 3074.25
$ws.Range("M61").Value = -2872.25
# Row 113
$ws.Range("H113").Value = 3074.25
$ws.Range("I113").Value = 3074.25
$ws.Range("K113").Value = 3074.25
$ws.Range("M113").Value = -904.25
# Row 119
$ws.Range("H119").Value = 85987
$ws.Range("J119").Value = 85987
$ws.Range("L119").Value = 85987
$ws.Range("N119").Value = -95663
# Row 128
$ws.Range("H128").Value = 66999.25
$ws.Range("J128").Value = 83998.5
$ws.Range("L128").Value = 83998.5
$ws.Range("N128").Value = -93958.5
# Row 141
$ws.Range("H141").Value = 99884.5
$ws.Range("J141").Value = 99884.5
$ws.Range("L141").Value = 99884.5
$ws.Range("N141").Value = -110244.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 51
$ws.Range("H51").Value = 30482.8
$ws.Range("I51").Value = 26485
$ws.Range("J51").Value = 33148
$ws.Range("K51").Value = 26485
$ws.Range("L51").Value = 33148
$ws.Range("M51").Value = -25975
$ws.Range("N51").Value = -34168
# Row 52
$ws.Range("H52").Value = 33979
$ws.Range("I52").Value = 11937
$ws.Range("K52").Value = 11937
$ws.Range("M52").Value = -11711
